# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rule "R40" (row 11) is renamed/re-keyed to "1" in the Rules sheet.
# The new value must stay a text string (not be auto-converted to a
# number), so the cell is formatted as Text before the value is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
